$d = $word.ActiveDocument

# Remove the stale "M2Doc version mismatch" warning block (including its
# surrounding spacer runs) that was left in the first paragraph from the
# outdated test template/runtime version check.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2
) | Out-Null
